# Update the "dSF" (column F) values for several rows as part of the
# "repull data, push all data, mean calculation" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F5").Value  = -1
$ws.Range("F6").Value  = -1
$ws.Range("F8").Value  = -1
$ws.Range("F9").Value  = 3
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 4
